$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The weekly data rows (2-4) get rotated: row2 <- row3, row3 <- row4, row4 <- row2
# for columns D (Fecha), J (Volumen), K (Precio minimo), L (Precio maximo),
# M (Precio promedio ponderado), P (Precio $/Kg).

$cols = @("D", "J", "K", "L", "M", "P")

# Capture original values first so we don't clobber data while rotating.
$row2 = @{}
$row3 = @{}
$row4 = @{}
foreach ($col in $cols) {
    $row2[$col] = $ws.Range($col + "2").Value2
    $row3[$col] = $ws.Range($col + "3").Value2
    $row4[$col] = $ws.Range($col + "4").Value2
}

foreach ($col in $cols) {
    $ws.Range($col + "2").Value2 = $row3[$col]
    $ws.Range($col + "3").Value2 = $row4[$col]
    $ws.Range($col + "4").Value2 = $row2[$col]
}
